$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the first "Dr. D. Y. Patil Pratishthan's" run into two runs with a
#    "_GoBack" bookmark in between (cursor-position artifact reproduced by
#    the commit): "Dr. D. Y. Patil P" | <bookmark _GoBack/> | "ratishthan's"
#    Scope the Find to the first paragraph only so the second, identical
#    heading further down the document is left untouched.
# ---------------------------------------------------------------------------
$titleRange = $d.Paragraphs.Item(1).Range
$titleFind = $titleRange.Find
$titleFound = $titleFind.Execute("Dr. D. Y. Patil P", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($titleFound) {
    $splitPoint = $d.Range($titleRange.End, $titleRange.End)
    $d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) "Himanshu Sharma & Pratik Warkari " -> "Himanshu Sharma"
#    (drop the " & Pratik Warkari " tail; the paragraph mark that used to
#    follow it is left in place, so the document keeps the same number of
#    paragraphs.)
# ---------------------------------------------------------------------------
$find1 = $d.Content.Find
$found1 = $find1.Execute(" & Pratik Warkari ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $find1.Parent.Text = ""
}

# ---------------------------------------------------------------------------
# 3) "Roll. No.: 13 & 32" -> "Roll. No.: 13"
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$found2 = $find2.Execute(" & 32", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $find2.Parent.Text = ""
}

# ---------------------------------------------------------------------------
# 4) Footer page-number field cached result: "4" -> "5"
#    Footer text lives in its own "story" (separate character numbering from
#    the main body), so address it through the footer Range's own
#    Characters collection rather than $d.Range(...).
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerRange = $footer.Range
if ($footerRange.Characters.Count -ge 1) {
    $pageNumChar = $footerRange.Characters.Item(1)
    if ($pageNumChar.Text -eq "4") {
        $pageNumChar.Text = "5"
    }
}
